$wb = $excel.ActiveWorkbook

# The "metadata" sheet holds the dataset_internal_id (B2) and
# indicator_internal_id (B3) values that need to be renamed.
$ws = $wb.Worksheets.Item("metadata")

$ws.Range("B2").Value = "LG"
$ws.Range("B3").Value = "LG.2F"
